$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared text used by every new run (results 005 - 011 per commit message)
$preprocess = 'convert unicode to ascii, remove break line, remove multiple spaces, convert to lower, space after punctuation, trim "space" and ","'
$features   = '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), #digit/#ascii, %kwName, %kwAddress, %kwPhone, #max_digit_skip_0 >= 7, first_character_ascii, first_character_digit, last_character_ascii, last_character_digit'
$model      = 'Neuron Network'
$modelDetails = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000'
$templateFilter = '0 filters: '

# Time, RunningTime(s), Test_Accuracy, Val_Accuracy, <unnamed col J>
$newRuns = @(
    @("20160415_171452", 1903.15,  0.995333333333333, 0.940594059405941, 0.475609756097561),
    @("20160415_174635", 1849.424, 0.996666666666667, 0.940594059405941, 0.597560975609756),
    @("20160415_181724", 1776.676, 0.995333333333333, 0.940594059405941, 0.451219512195122),
    @("20160415_184701", 1331.827, 0.994666666666667, 0.940594059405941, 0.451219512195122),
    @("20160415_190913", 1246.841, 0.988,              0.940594059405941, 0.48780487804878)
)

$startRow = 17
for ($i = 0; $i -lt $newRuns.Length; $i++) {
    $row = $startRow + $i
    $run = $newRuns[$i]

    $ws.Cells.Item($row, 1).Value  = $run[0]        # Time
    $ws.Cells.Item($row, 2).Value  = $run[1]        # RunningTime(s)
    $ws.Cells.Item($row, 3).Value  = $preprocess    # Preprocess
    $ws.Cells.Item($row, 4).Value  = $features      # Features
    $ws.Cells.Item($row, 5).Value  = $model         # Model
    $ws.Cells.Item($row, 6).Value  = $modelDetails  # Model_Details
    $ws.Cells.Item($row, 7).Value  = $run[2]        # Test_Accuracy
    $ws.Cells.Item($row, 8).Value  = $run[3]        # Val_Accuracy
    $ws.Cells.Item($row, 9).Value  = $templateFilter# Template Filter
    $ws.Cells.Item($row, 10).Value = $run[4]        # (unnamed column J)
}
